$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 26 with the LeetCode problem "Peak Index in a Mountain Array"
$ws.Cells.Item(26, 1).Value = 852
$ws.Cells.Item(26, 2).Value = "Peak Index in a Mountain Array"
$ws.Cells.Item(26, 3).Value = "#array #binary-search #重点 "
$ws.Cells.Item(26, 4).Value = "medium"
$ws.Cells.Item(26, 5).Value = 3
$ws.Cells.Item(26, 6).Value = 0
$ws.Cells.Item(26, 7).Value = 3
$ws.Cells.Item(26, 8).Value = 45838
$ws.Cells.Item(26, 9).Value = 45838

# Match the style/format used by the row above (row 25)
$ws.Range("A26").Style = $ws.Range("A25").Style
$ws.Range("B26").Style = $ws.Range("B25").Style
$ws.Range("C26").Style = $ws.Range("C25").Style
$ws.Range("D26").Style = $ws.Range("D25").Style
$ws.Range("E26").Style = $ws.Range("E25").Style
$ws.Range("F26").Style = $ws.Range("F25").Style
$ws.Range("G26").Style = $ws.Range("G25").Style
$ws.Range("H26").Style = $ws.Range("H25").Style
$ws.Range("I26").Style = $ws.Range("I25").Style

$ws.Range("H26:I26").NumberFormat = $ws.Range("H25:I25").NumberFormat

$ws.Rows.Item(26).RowHeight = $ws.Rows.Item(25).RowHeight

# Update the active selection to reflect the new state
$ws.Range("D30").Select()
